# Add season record columns (Wins, Losses, Ties) for each player/team row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD, AE, AF -----------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the existing header cells (bold, bordered, centered/top)
# by copying the formatting from an existing header cell.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows (2-49): season record values for every player ---------------
$lastRow = 49
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 90   # column AD -> Wins
    $ws.Cells.Item($r, 31).Value = 72   # column AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # column AF -> Ties
}
